# Auto-generated script to apply scheduled market-data refresh values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 499.4
$ws.Range("I12").Value = 499.4
$ws.Range("K12").Value = 499.4
$ws.Range("M12").Value = -329.4
$ws.Range("H33").Value = 400.85715
$ws.Range("I33").Value = 400.85715
$ws.Range("K33").Value = 400.85715
$ws.Range("M33").Value = -171.85715
$ws.Range("H132").Value = 4856.8696
$ws.Range("I132").Value = 4804.909
$ws.Range("K132").Value = 14414.727
$ws.Range("M132").Value = -11884.727
$ws.Range("H137").Value = 3229.55
$ws.Range("I137").Value = 2593
$ws.Range("J137").Value = 5139.2
$ws.Range("K137").Value = 7779
$ws.Range("L137").Value = 15417.6
$ws.Range("M137").Value = -5229
$ws.Range("N137").Value = -20517.6
$ws.Range("H138").Value = 2525.3684
$ws.Range("I138").Value = 1876
$ws.Range("J138").Value = 2997.6365
$ws.Range("K138").Value = 5628
$ws.Range("L138").Value = 8992.9095
$ws.Range("M138").Value = -488
$ws.Range("N138").Value = -19272.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1288.1724
$ws.Range("I32").Value = 1314.4814
$ws.Range("K32").Value = 1314.4814
$ws.Range("M32").Value = -1027.4814
$ws.Range("H122").Value = 1834.875
$ws.Range("I122").Value = 1834.875
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5504.625
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3054.625
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 72333
$ws.Range("J58").Value = 72333
$ws.Range("L58").Value = 72333
$ws.Range("N58").Value = -72921
$ws.Range("H60").Value = 26000
$ws.Range("I60").Value = 26000
$ws.Range("K60").Value = 26000
$ws.Range("M60").Value = -25401
$ws.Range("H74").Value = 80590
$ws.Range("J74").Value = 80590
$ws.Range("L74").Value = 80590
$ws.Range("N74").Value = -82462
$ws.Range("H77").Value = 80590
$ws.Range("J77").Value = 80590
$ws.Range("L77").Value = 241770
$ws.Range("N77").Value = -251130
$ws.Range("H92").Value = 49687.625
$ws.Range("J92").Value = 49687.625
$ws.Range("L92").Value = 49687.625
$ws.Range("N92").Value = -54679.625
$ws.Range("H94").Value = 1210
$ws.Range("I94").Value = 1236.25
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 1236.25
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -785.25
$ws.Range("N94").Value = -1902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2332.1667
$ws.Range("J31").Value = 998.75
$ws.Range("L31").Value = 998.75
$ws.Range("N31").Value = -1588.75
$ws.Range("H34").Value = 2332.1667
$ws.Range("J34").Value = 998.75
$ws.Range("L34").Value = 998.75
$ws.Range("N34").Value = -1402.75
$ws.Range("H41").Value = 27323.4
$ws.Range("J41").Value = 69999
$ws.Range("L41").Value = 69999
$ws.Range("N41").Value = -70855
$ws.Range("H47").Value = 8544.5
$ws.Range("I47").Value = 8392.666999999999
$ws.Range("K47").Value = 8392.666999999999
$ws.Range("M47").Value = -7826.666999999999
$ws.Range("H58").Value = 3581.5
$ws.Range("I58").Value = 2872.25
$ws.Range("K58").Value = 2872.25
$ws.Range("M58").Value = -2669.25
$ws.Range("H60").Value = 42228.145
$ws.Range("J60").Value = 59899.25
$ws.Range("L60").Value = 59899.25
$ws.Range("N60").Value = -60921.25
$ws.Range("H92").Value = 20201
$ws.Range("J92").Value = 20201
$ws.Range("L92").Value = 20201
$ws.Range("N92").Value = -25193
$ws.Range("H108").Value = 70000
$ws.Range("J108").Value = 70000
$ws.Range("L108").Value = 70000
$ws.Range("N108").Value = -77680
$ws.Range("H111").Value = 33777.5
$ws.Range("J111").Value = 33777.5
$ws.Range("L111").Value = 33777.5
$ws.Range("N111").Value = -41957.5
$ws.Range("H134").Value = 4252.8823
$ws.Range("I134").Value = 4081.1875
$ws.Range("K134").Value = 12243.5625
$ws.Range("M134").Value = -9708.5625
$ws.Range("H136").Value = 3581.5
$ws.Range("I136").Value = 2872.25
$ws.Range("K136").Value = 8616.75
$ws.Range("M136").Value = -6066.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 882
$ws.Range("I11").Value = 882
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2646
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -2506
$ws.Range("N11").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7207309.5
$ws.Range("I11").Value = 9603046
$ws.Range("J11").Value = 20099.4
$ws.Range("K11").Value = 9603046
$ws.Range("L11").Value = 20099.4
$ws.Range("M11").Value = -9602907
$ws.Range("N11").Value = -20377.4
$ws.Range("H21").Value = 49949.25
$ws.Range("J21").Value = 49949.25
$ws.Range("L21").Value = 49949.25
$ws.Range("N21").Value = -50295.25
$ws.Range("H24").Value = 22433.232
$ws.Range("J24").Value = 25590.75
$ws.Range("L24").Value = 25590.75
$ws.Range("N24").Value = -25936.75
$ws.Range("H30").Value = 49949.25
$ws.Range("J30").Value = 49949.25
$ws.Range("L30").Value = 49949.25
$ws.Range("N30").Value = -50159.25
$ws.Range("H47").Value = 57500
$ws.Range("I47").Value = 49000
$ws.Range("J47").Value = 100000
$ws.Range("K47").Value = 49000
$ws.Range("L47").Value = 100000
$ws.Range("M47").Value = -48432
$ws.Range("N47").Value = -101136
$ws.Range("H92").Value = 10901
$ws.Range("J92").Value = 11044.286
$ws.Range("L92").Value = 11044.286
$ws.Range("N92").Value = -14788.286
$ws.Range("H97").Value = 380.42856
$ws.Range("I97").Value = 296.55554
$ws.Range("J97").Value = 531.4
$ws.Range("K97").Value = 296.55554
$ws.Range("L97").Value = 531.4
$ws.Range("M97").Value = 199.44446
$ws.Range("N97").Value = -1523.4
$ws.Range("H126").Value = 2234.4443
$ws.Range("I126").Value = 2262
$ws.Range("J126").Value = 2014
$ws.Range("K126").Value = 6786
$ws.Range("L126").Value = 6042
$ws.Range("M126").Value = -4316
$ws.Range("N126").Value = -10982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 88999.5
$ws.Range("J36").Value = 88999.5
$ws.Range("L36").Value = 88999.5
$ws.Range("N36").Value = -90123.5
$ws.Range("H40").Value = 14984.823
$ws.Range("I40").Value = 14365.223
$ws.Range("J40").Value = 15681.875
$ws.Range("K40").Value = 14365.223
$ws.Range("L40").Value = 15681.875
$ws.Range("M40").Value = -14229.223
$ws.Range("N40").Value = -15953.875
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2000
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 10000
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -17488
$ws.Range("H100").Value = 1936.6666
$ws.Range("I100").Value = 1936.6666
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1936.6666
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1395.6666
$ws.Range("N100").ClearContents()
$ws.Range("H103").Value = 8918.4
$ws.Range("J103").Value = 8918.4
$ws.Range("L103").Value = 8918.4
$ws.Range("N103").Value = -11262.4
$ws.Range("H139").Value = 65000
$ws.Range("J139").Value = 65000
$ws.Range("L139").Value = 65000
$ws.Range("N139").Value = -75280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 27959.6
$ws.Range("I62").Value = 9999
$ws.Range("J62").Value = 32449.75
$ws.Range("K62").Value = 9999
$ws.Range("L62").Value = 32449.75
$ws.Range("M62").Value = -9375
$ws.Range("N62").Value = -33697.75
$ws.Range("H65").Value = 27959.6
$ws.Range("I65").Value = 9999
$ws.Range("J65").Value = 32449.75
$ws.Range("K65").Value = 49995
$ws.Range("L65").Value = 162248.75
$ws.Range("M65").Value = -46875
$ws.Range("N65").Value = -168488.75
$ws.Range("H104").Value = 22974.834
$ws.Range("J104").Value = 22974.834
$ws.Range("L104").Value = 22974.834
$ws.Range("N104").Value = -29962.834
$ws.Range("H110").Value = 49999.5
$ws.Range("J110").Value = 49999.5
$ws.Range("L110").Value = 49999.5
$ws.Range("N110").Value = -58179.5
$ws.Range("H132").Value = 8332.333000000001
$ws.Range("I132").Value = 6197.8
$ws.Range("K132").Value = 18593.4
$ws.Range("M132").Value = -16063.4
$ws.Range("H136").Value = 2999.0322
$ws.Range("I136").Value = 3142.2307
$ws.Range("K136").Value = 9426.6921
$ws.Range("M136").Value = -6876.6921
